# Daily attendance processing
# Normalizes the "Recorded By" (column G) lists on the active sheet by
# sorting the comma-separated recorder names/emails into ordinal
# (case-sensitive, ASCII) ascending order, e.g.
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#   "backup@backdoor.com, system, System" -> "System, backup@backdoor.com, system"

function Compare-Ordinal($s1, $s2) {
    $len1 = $s1.Length
    $len2 = $s2.Length
    $minLen = [Math]::Min($len1, $len2)
    for ($ci = 0; $ci -lt $minLen; $ci++) {
        $c1 = [int][char]$s1[$ci]
        $c2 = [int][char]$s2[$ci]
        if ($c1 -lt $c2) { return -1 }
        if ($c1 -gt $c2) { return 1 }
    }
    if ($len1 -lt $len2) { return -1 }
    if ($len1 -gt $len2) { return 1 }
    return 0
}

function Sort-Ordinal($arr) {
    $result = @($arr)
    $n = $result.Count
    for ($oi = 1; $oi -lt $n; $oi++) {
        $key = $result[$oi]
        $j = $oi - 1
        while ($j -ge 0 -and (Compare-Ordinal $result[$j] $key) -gt 0) {
            $result[$j + 1] = $result[$j]
            $j = $j - 1
        }
        $result[$j + 1] = $key
    }
    return $result
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $sortedParts = Sort-Ordinal $parts
            $newVal = [string]::Join(", ", $sortedParts)
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}

Write-Output "Recorded By column normalized for rows 2..$lastRow"
